$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name / link text updates (rows that were reordered)
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('B34').Value = 'Dai'
$ws.Range('C34').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'

# Price / volume(1h) updates (numeric-looking text, protected with a leading
# apostrophe so Excel keeps them as text, then style reset to drop the quote-prefix flag)
$ws.Range('D2').Value = "'61.551.49"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -1.04%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.378.48"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -1.67%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.11%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'407.18"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -1.79%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'136.68"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +11.12%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.594"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.57%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  -0.09%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.673"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +2.60%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.122"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -5.07%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'42.97"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +4.06%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  -1.17%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'3.904.65"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -2.19%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'8.40"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -1.26%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'19.73"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -0.10%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'3.391.42"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -1.39%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'61.479.14"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -1.30%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.72%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'11.01"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +1.32%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  -2.99%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  -3.04%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'83.92"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.57%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'315.57"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -1.55%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'12.82"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -1.21%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  -0.75%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  +11.66%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'8.33"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +5.60%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  -5.30%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'7.71"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -2.25%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'0.172"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -0.53%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'0.117"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +0.28%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'2.59"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.94%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'11.36"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -1.24%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'1.00"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.02%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'41.01"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -2.49%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.0481"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.63%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'51.94"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -1.48%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'0.998"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -0.26%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  -2.06%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  -2.86%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'138.27"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +2.79%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'1.97"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -0.32%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.124"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -0.86%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.296"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +4.51%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'4.05"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +4.07%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'16.74"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -2.44%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  +0.94%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'21.55"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -2.70%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'2.131.54"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -3.15%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -5.06%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'1.93"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +1.77%  "
$ws.Range('E51').Style = 'Normal'
